$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set C1 header style to match B1 (bold, centered)
$ws.Range("C1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("C1").Font.Bold = $true

# Set the created date for rows 2-8 in column C, formatted as a date (m/d/yyyy),
# centered to match the header style, and leave D (Last Edited) blank but
# pre-formatted as a date as well.
$createdDate = Get-Date -Year 2020 -Month 9 -Day 1 -Hour 0 -Minute 0 -Second 0

for ($r = 2; $r -le 8; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Value = $createdDate
    $cCell.NumberFormat = "m/d/yyyy"
    $cCell.HorizontalAlignment = -4108  # xlCenter

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "m/d/yyyy"
}

# Update active selection to C1
$ws.Range("C1").Select()
